$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.Value = "'" + $value
    $cell.Style = $ws.Range("D4").Style
}

Set-TextValue $ws.Range("D2") "67.537.83"
Set-TextValue $ws.Range("E2") "  +0.08%  "
Set-TextValue $ws.Range("D3") "3.529.40"
Set-TextValue $ws.Range("E3") "  +0.88%  "
Set-TextValue $ws.Range("E4") "  +0.44%  "
Set-TextValue $ws.Range("D5") "613.73"
Set-TextValue $ws.Range("E5") "  +0.68%  "
Set-TextValue $ws.Range("D6") "152.45"
Set-TextValue $ws.Range("E6") "  +2.14%  "
Set-TextValue $ws.Range("D7") "3.520.68"
Set-TextValue $ws.Range("E7") "  +0.63%  "
Set-TextValue $ws.Range("E8") "  +0.28%  "
Set-TextValue $ws.Range("D9") "0.485"
Set-TextValue $ws.Range("E9") "  +0.99%  "
Set-TextValue $ws.Range("D10") "0.144"
Set-TextValue $ws.Range("E10") "  +4.21%  "
Set-TextValue $ws.Range("D11") "7.18"
Set-TextValue $ws.Range("E11") "  +3.55%  "
Set-TextValue $ws.Range("D12") "0.434"
Set-TextValue $ws.Range("E12") "  +2.65%  "
Set-TextValue $ws.Range("D13") "0.0000220"
Set-TextValue $ws.Range("E13") "  +1.21%  "
Set-TextValue $ws.Range("D14") "32.68"
Set-TextValue $ws.Range("E14") "  +3.55%  "
Set-TextValue $ws.Range("D15") "4.134.41"
Set-TextValue $ws.Range("E15") "  +1.11%  "
Set-TextValue $ws.Range("D16") "3.553.64"
Set-TextValue $ws.Range("E16") "  +1.68%  "
Set-TextValue $ws.Range("D17") "68.302.38"
Set-TextValue $ws.Range("E17") "  +1.46%  "
Set-TextValue $ws.Range("D18") "0.116"
Set-TextValue $ws.Range("E18") "  -0.28%  "
Set-TextValue $ws.Range("D19") "6.68"
Set-TextValue $ws.Range("E19") "  +3.91%  "
Set-TextValue $ws.Range("D20") "15.71"
Set-TextValue $ws.Range("E20") "  +4.47%  "
Set-TextValue $ws.Range("D21") "9.85"
Set-TextValue $ws.Range("E21") "  +9.00%  "
Set-TextValue $ws.Range("D22") "450.91"
Set-TextValue $ws.Range("E22") "  +0.71%  "
Set-TextValue $ws.Range("D23") "0.634"
Set-TextValue $ws.Range("E23") "  +1.53%  "
Set-TextValue $ws.Range("D24") "77.85"
Set-TextValue $ws.Range("E24") "  +0.77%  "
Set-TextValue $ws.Range("B25") "WrappedeETH"
Set-TextValue $ws.Range("C25") "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue $ws.Range("D25") "3.688.95"
Set-TextValue $ws.Range("E25") "  +1.44%  "
Set-TextValue $ws.Range("B26") "PEPE"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D26") "0.0000129"
Set-TextValue $ws.Range("E26") "  -0.18%  "
Set-TextValue $ws.Range("D27") "10.37"
Set-TextValue $ws.Range("E27") "  +2.02%  "
Set-TextValue $ws.Range("E28") "  -0.24%  "
Set-TextValue $ws.Range("D29") "8.92"
Set-TextValue $ws.Range("E29") "  +6.48%  "
Set-TextValue $ws.Range("D30") "2.55"
Set-TextValue $ws.Range("E30") "  +1.52%  "
Set-TextValue $ws.Range("D31") "1.66"
Set-TextValue $ws.Range("E31") "  +4.66%  "
Set-TextValue $ws.Range("D32") "0.168"
Set-TextValue $ws.Range("E32") "  +1.63%  "
Set-TextValue $ws.Range("E33") "  +0.43%  "
Set-TextValue $ws.Range("D34") "6.34"
Set-TextValue $ws.Range("E34") "  +3.43%  "
Set-TextValue $ws.Range("D35") "25.83"
Set-TextValue $ws.Range("E35") "  +0.10%  "
Set-TextValue $ws.Range("D36") "1.90"
Set-TextValue $ws.Range("E36") "  +1.82%  "
Set-TextValue $ws.Range("D37") "3.527.42"
Set-TextValue $ws.Range("E37") "  +1.06%  "
Set-TextValue $ws.Range("D38") "8.16"
Set-TextValue $ws.Range("E38") "  +1.83%  "
Set-TextValue $ws.Range("E39") "  +0.03%  "
Set-TextValue $ws.Range("D40") "2.32"
Set-TextValue $ws.Range("E40") "  +4.67%  "
Set-TextValue $ws.Range("B41") "FirstDigitalUSD"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D41") "1.00"
Set-TextValue $ws.Range("E41") "  +0.60%  "
Set-TextValue $ws.Range("B42") "Monero"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D42") "177.56"
Set-TextValue $ws.Range("E42") "  +4.57%  "
Set-TextValue $ws.Range("D43") "0.0908"
Set-TextValue $ws.Range("E43") "  +3.48%  "
Set-TextValue $ws.Range("D44") "5.49"
Set-TextValue $ws.Range("E44") "  +0.99%  "
Set-TextValue $ws.Range("B45") "InjectiveProtocol"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D45") "30.22"
Set-TextValue $ws.Range("E45") "  +13.03%  "
Set-TextValue $ws.Range("B46") "Mantle"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D46") "0.886"
Set-TextValue $ws.Range("E46") "  +0.07%  "
Set-TextValue $ws.Range("B47") "OKB"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D47") "45.36"
Set-TextValue $ws.Range("E47") "  -0.37%  "
Set-TextValue $ws.Range("B48") "ONDO"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws.Range("D48") "1.32"
Set-TextValue $ws.Range("E48") "  +4.64%  "
Set-TextValue $ws.Range("D49") "2.63"
Set-TextValue $ws.Range("E49") "  +2.61%  "
Set-TextValue $ws.Range("D50") "7.71"
Set-TextValue $ws.Range("E50") "  +2.23%  "
Set-TextValue $ws.Range("D51") "0.258"
Set-TextValue $ws.Range("E51") "  +5.70%  "
